$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2862
$ws.Range("F3").Value = 1159
$ws.Range("F4").Value = 20913
$ws.Range("F5").Value = 99
$ws.Range("F6").Value = 2736
$ws.Range("F7").Value = 795
$ws.Range("F9").Value = 502
$ws.Range("F10").Value = 753
$ws.Range("F11").Value = 275
$ws.Range("F13").Value = 71
$ws.Range("F15").Value = 507
$ws.Range("F17").Value = 254
$ws.Range("F18").Value = 12
$ws.Range("F19").Value = 415
$ws.Range("F20").Value = 41
$ws.Range("F22").Value = 27
$ws.Range("F23").Value = 20

# Sheet "演出" (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 327
$ws.Range("F12").Value = 96
$ws.Range("F14").Value = 141

# Sheet "本地生活" (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6114
$ws.Range("F3").Value = 693
$ws.Range("F4").Value = 673
$ws.Range("F5").Value = 1531
$ws.Range("F6").Value = 51

# Sheet "全部类型" (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6114
$ws.Range("F3").Value = 693
$ws.Range("F4").Value = 673
$ws.Range("F5").Value = 1531
$ws.Range("F6").Value = 2862
$ws.Range("F7").Value = 1159
$ws.Range("F8").Value = 20913
$ws.Range("F11").Value = 99
$ws.Range("F13").Value = 327
$ws.Range("F14").Value = 2736
$ws.Range("F15").Value = 795
$ws.Range("F17").Value = 51
$ws.Range("F19").Value = 502
$ws.Range("F20").Value = 753
$ws.Range("F21").Value = 275
$ws.Range("F24").Value = 71
$ws.Range("F30").Value = 507
$ws.Range("F31").Value = 96
$ws.Range("F34").Value = 254
$ws.Range("F35").Value = 141
$ws.Range("F36").Value = 141
$ws.Range("F37").Value = 12
$ws.Range("F38").Value = 415
$ws.Range("F43").Value = 27
$ws.Range("F44").Value = 20
